$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account holder info ---
$ws.Range("C2").Value = "Hartmut"

# B3 is a purely-numeric-looking string ("2570314725427075"). Assigning it
# directly via .Value would get auto-coerced into a real number by Excel's
# General number format. Route it through TEXT()+paste-values instead so it
# lands back down as a genuine text cell (matching the original inline
# string type) while keeping the cell's existing style index intact.
$ws.Range("B3").Formula = '=TEXT("2570314725427075","@")'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 17.02.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "21.02."
$ws.Range("C6").Value = "22.02."
$ws.Range("D6").Value = "BURGER KING Soest"
$ws.Range("E6").Value = "21,76-"

# --- Row 7 ---
$ws.Range("B7").Value = "25.02."
$ws.Range("C7").Value = "26.02."
$ws.Range("D7").Value = "PAYPAL AYBKJN"
$ws.Range("E7").Value = "37,44-"

# --- Row 8 ---
$ws.Range("B8").Value = "28.02."
$ws.Range("C8").Value = "29.02."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU CEQMRW"
$ws.Range("E8").Value = "150,34-"

# --- Row 9 ---
$ws.Range("B9").Value = "29.02."
$ws.Range("C9").Value = "01.03."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-55073539"
$ws.Range("E9").Value = "57,04-"

# --- Rows 10 & 11 lost their transactions entirely; cells go blank. ---
# B/C/D keep their original style (8); E switches from the plain
# right-aligned style (17) to the right+vcenter+wrap style (12), matching
# the other blank amount cells in the sheet.
$ws.Range("B10:D11").Value = ""
$rngE = $ws.Range("E10:E11")
$rngE.Value = ""
$rngE.WrapText = $true
$rngE.VerticalAlignment = -4108

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 05.03.2024"
$ws.Range("E12").Value = "266,58-"

# --- Next billing date note ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.03.2024"
